# Add a new "Team Meetings" sheet (team discussion minutes) as the first
# sheet in the workbook, ahead of Provencher / Philippona / Ningge.
#
# The easiest way to get a sheet whose column widths / date+hours number
# formats / header font match the other timesheets is to clone the
# "Provencher" sheet (same 3-column DATE / Time Spent / What I worked on
# layout) and then overwrite its data.

$wb = $excel.ActiveWorkbook

$provencher = $wb.Worksheets.Item("Provencher")

# Clone Provencher -> new sheet placed immediately before it. Excel numbers
# the clone's sheetId 4 (next unused id) and slots it in as rId1, pushing
# Provencher/Philippona/Ningge to rId2/rId3/rId4 - exactly what we want.
$provencher.Copy($provencher)
$meetings = $wb.Worksheets.Item(1)
$meetings.Name = "Team Meetings"

# Make room for a new row 1 ("Team Discussions" banner) above the existing
# DATE / Time Spent / What I worked on header row, then trim the sheet back
# down to 21 rows total (2 header rows + 12 data rows + 7 blank rows) to
# match the new, shorter table.
$meetings.Rows("1").Insert()
$meetings.Rows("22:39").Delete()

# A1: bold "Team Discussions" banner - reuse the bold header style already
# sitting on row 2 (col A) rather than re-deriving a font style by hand.
$meetings.Range("A2").Copy()
$meetings.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$meetings.Range("A1").Value = "Team Discussions"

# Fill in the "What I worked on" text column first, in the same order the
# original author's workbook lists its shared strings (rows 3-7, 9-11,
# 13-14, then 8, then 12 - the two meeting notes that were recorded out of
# sequence), so new shared-string entries land at the expected indices.
$meetings.Range("C3").Value = "Initial project discussion. Idea for proceduraly generated game. 3D RPG"
$meetings.Range("C4").Value = "Scrapped ambitions for 3D game, decided on 2D procedural generation."
$meetings.Range("C5").Value = "Discussed story, art and inspirations for game. "
$meetings.Range("C6").Value = "Game Proposal Presentation"
$meetings.Range("C7").Value = "Discussion on task splitting and implementation details."
$meetings.Range("C9").Value = "Online discussion on map creation."
$meetings.Range("C10").Value = "Online discussion on player control"
$meetings.Range("C11").Value = "Online discussion on enemy types"
$meetings.Range("C13").Value = "Online collaborative trouble shooting on various issues."
$meetings.Range("C14").Value = "Discussion on remaining tasks and issues"
$meetings.Range("C8").Value = "Setting up of online repository and issue tracker establishment"
$meetings.Range("C12").Value = "Playtesting and discussion of glitches and mechanic improvement ideas"

# DATE column.
$meetings.Range("A3").Value = 42262
$meetings.Range("A4").Value = 42265
$meetings.Range("A5").Value = 42267
$meetings.Range("A6").Value = 42269
$meetings.Range("A7").Value = 42276
$meetings.Range("A8").Value = 42291
$meetings.Range("A9").Value = 42292
$meetings.Range("A10").Value = 42295
$meetings.Range("A11").Value = 42297
$meetings.Range("A12").Value = 42300
$meetings.Range("A13").Value = 42301
$meetings.Range("A14").Value = 42304

# Time Spent (hours) column.
$meetings.Range("B3").Value = 0.5
$meetings.Range("B4").Value = 0.25
$meetings.Range("B5").Value = 1
$meetings.Range("B6").Value = 0.25
$meetings.Range("B7").Value = 0.5
$meetings.Range("B8").Value = 1.5
$meetings.Range("B9").Value = 0.5
$meetings.Range("B10").Value = 0.25
$meetings.Range("B11").Value = 0.25
$meetings.Range("B12").Value = 1.5
$meetings.Range("B13").Value = 1
$meetings.Range("B14").Value = 0.5

# Rows 3-13 show hours to two decimal places; row 14 keeps the one-decimal
# format inherited from the cloned Provencher sheet.
$meetings.Range("B3:B13").NumberFormat = "0.00"

# Restore each sheet's own remembered selection: Provencher keeps its own
# last selected range, and Team Meetings (now the active tab) ends up with
# C16 selected.
$provencher.Range("A1:C18").Select()
$meetings.Activate()
$meetings.Range("C16").Select()
